$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.317.47'
$ws.Range("E2").Value = '  +0.89%  '

$ws.Range("D3").Value = '1.665.41'
$ws.Range("E3").Value = '  +0.73%  '

$ws.Range("E4").Value = '  +0.82%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '219.23'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.81%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5350'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.87%  '

$ws.Range("E7").Value = '  +0.76%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2662'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.43%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06411'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.19%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.64'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.36%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07834'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.49%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.570'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.57%  '

$ws.Range("D13").Value = '1.673.57'
$ws.Range("E13").Value = '  +1.22%  '

$ws.Range("D14").Value = '1.892.54'
$ws.Range("E14").Value = '  +0.62%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5530'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.01%  '

$ws.Range("D16").Value = '0.0₅8237'
$ws.Range("E16").Value = '  +0.01%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.80'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.55%  '

$ws.Range("E18").Value = '  +0.81%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.700'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.81%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '193.92'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.62%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.26'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.03%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.042'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.31%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.011'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.76%  '

$ws.Range("E24").Value = '  +3.30%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1230'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.01%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.198'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.48%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.12'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.41%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.483'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.89%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.05837'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.45%  '

$ws.Range("E30").Value = '  +0.90%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.626'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.65%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.283'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.72%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.618'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.46%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9688'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.13%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.825'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.63%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.419'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.29%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5819'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.71%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01606'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.60%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.8693'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.92%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.875'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.68%  '

$ws.Range("D41").Value = '1.053.61'
$ws.Range("E41").Value = '  +2.66%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '105.16'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.91%  '

$ws.Range("E43").Value = '  +0.83%  '

$ws.Range("D44").Value = '1.803.94'
$ws.Range("E44").Value = '  +0.42%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '57.89'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.46%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.015'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.36%  '

$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4386'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.69%  '

$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.030'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.76%  '

$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = '0.0₈103'
$ws.Range("E49").Value = '  -8.34%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05164'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.34%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.412'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.70%  '
